# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for a set of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12  = @("aa", "Agree/Accept")
    15  = @("sd", "Statement-non-opinion")
    20  = @("aa", "Agree/Accept")
    23  = @("ba", "Appreciation")
    30  = @("sd", "Statement-non-opinion")
    35  = @("sd", "Statement-non-opinion")
    37  = @("sd", "Statement-non-opinion")
    41  = @("%", "Uninterpretable")
    69  = @("sv", "Statement-opinion")
    70  = @("sd", "Statement-non-opinion")
    73  = @("sd", "Statement-non-opinion")
    75  = @("%", "Uninterpretable")
    79  = @("b", "Acknowledge (Backchannel)")
    91  = @("aa", "Agree/Accept")
    106 = @("aa", "Agree/Accept")
    108 = @("sd", "Statement-non-opinion")
    118 = @("sd", "Statement-non-opinion")
    120 = @("qy", "Yes-No-Question")
    121 = @("sd", "Statement-non-opinion")
    134 = @("sd", "Statement-non-opinion")
    135 = @("sv", "Statement-opinion")
    137 = @("sd", "Statement-non-opinion")
    140 = @("sd", "Statement-non-opinion")
    145 = @("sv", "Statement-opinion")
    147 = @("b", "Acknowledge (Backchannel)")
    148 = @("sd", "Statement-non-opinion")
    150 = @("sd", "Statement-non-opinion")
    160 = @("sv", "Statement-opinion")
    167 = @("%", "Uninterpretable")
    170 = @("sv", "Statement-opinion")
    175 = @("%", "Uninterpretable")
    176 = @("sd", "Statement-non-opinion")
    183 = @("aa", "Agree/Accept")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("I$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
}
